$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $style = $p.Range.ParagraphStyle.NameLocal
    if ($style -eq "Author" -or $style -eq "Date") {
        $p.Range.Delete()
    }
}
